# Updated Master data as per 16th May Refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of test data (rows 34-36), continuing the
# pattern already present in the sheet for regcntr_id = 10005.
$newRows = @(
    @(10005, 110033, 10005, "eng", $true, "superadmin", "now()"),
    @(10005, 110034, 10005, "eng", $true, "superadmin", "now()"),
    @(10005, 110035, 10005, "eng", $true, "superadmin", "now()")
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Move the active selection below the newly pasted data, matching the
# post-paste selection state (full rows selected starting at row 37).
$ws.Range("A37:XFD1048576").Select()
